$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 11 with the new "n dernières mesures" route. Order matters for
# the shared-strings table layout: description, then example URL, then the
# route URL itself (matches how Excel appended new unique strings on save).
$ws.Range("D11").Value = "Recuperer les {n} dernières mesures classés chronologiquement (de la plus ancienne à la plus récente mesure) pour la plante ayant l'adresse Mac {Adresse_Mac_Plante}"
$ws.Range("E11").Value = "https://azammouri.com/pc/uploads/mesures.php?id=AC:67:B2:36:61:D4&nbrMesure=20"
$ws.Range("A11").Value = "https://azammouri.com/pc/uploads/plante.php?id={Adresse_Mac_Plante}&nbrMesure={n}"
$ws.Range("B11").Value = "GET"
$ws.Range("C11").Value = "Json"

# Turn A11 into a hyperlink pointing at the same URL, matching the style of the
# other "route" cells in column A (their link addresses are percent-encoded,
# e.g. {Adresse_Mac_Plante} -> %7bAdresse_Mac_Plante%7d, same as A9/A10).
# Adding the hyperlink stamps the built-in "Lien hypertexte" font as direct
# formatting, so re-apply the plain vertical-center style used by the other
# route cells (e.g. A9) afterwards.
$ws.Hyperlinks.Add($ws.Range("A11"), "https://azammouri.com/pc/uploads/plante.php?id=%7bAdresse_Mac_Plante%7d&nbrMesure=%7bn%7d")
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Update the selection to match the saved view state after the edit (also
# drops the stale topLeftCell="A4" scroll position from the old selection).
$ws.Range("B14").Select()
